# Update "想去人数" (want-to-go headcount) figures on the 展览 (Exhibition)
# and 全部类型 (All types) sheets to the refreshed crawl numbers.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "展览";     Row = 6;  Value = 944 },
    @{ Sheet = "展览";     Row = 8;  Value = 305 },
    @{ Sheet = "展览";     Row = 9;  Value = 1232 },
    @{ Sheet = "展览";     Row = 11; Value = 281 },
    @{ Sheet = "展览";     Row = 18; Value = 7744 },
    @{ Sheet = "展览";     Row = 21; Value = 4405 },
    @{ Sheet = "展览";     Row = 23; Value = 2247 },
    @{ Sheet = "展览";     Row = 24; Value = 954 },
    @{ Sheet = "展览";     Row = 26; Value = 236 },
    @{ Sheet = "展览";     Row = 33; Value = 1926 },
    @{ Sheet = "展览";     Row = 41; Value = 2036 },
    @{ Sheet = "全部类型"; Row = 8;  Value = 944 },
    @{ Sheet = "全部类型"; Row = 10; Value = 305 },
    @{ Sheet = "全部类型"; Row = 11; Value = 1232 },
    @{ Sheet = "全部类型"; Row = 13; Value = 281 },
    @{ Sheet = "全部类型"; Row = 20; Value = 7744 },
    @{ Sheet = "全部类型"; Row = 23; Value = 4405 },
    @{ Sheet = "全部类型"; Row = 25; Value = 2247 },
    @{ Sheet = "全部类型"; Row = 26; Value = 954 },
    @{ Sheet = "全部类型"; Row = 28; Value = 236 },
    @{ Sheet = "全部类型"; Row = 36; Value = 1926 },
    @{ Sheet = "全部类型"; Row = 45; Value = 2036 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Cells.Item($u.Row, 6).Value = $u.Value
}
